$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10068.923
$ws.Range("I64").Value = 7223.25
$ws.Range("K64").Value = 7223.25
$ws.Range("M64").Value = -6975.25

$ws.Range("H67").Value = 10068.923
$ws.Range("I67").Value = 7223.25
$ws.Range("K67").Value = 7223.25
$ws.Range("M67").Value = -6365.25

$ws.Range("H92").Value = 1165.5641
$ws.Range("I92").Value = 615.2
$ws.Range("J92").Value = 3000.111
$ws.Range("K92").Value = 615.2
$ws.Range("L92").Value = 3000.111
$ws.Range("M92").Value = 632.8
$ws.Range("N92").Value = -5496.111

$ws.Range("H100").Value = 3100.84
$ws.Range("I100").Value = 2751.8572
$ws.Range("J100").Value = 3545
$ws.Range("K100").Value = 2751.8572
$ws.Range("L100").Value = 3545
$ws.Range("M100").Value = -2210.8572
$ws.Range("N100").Value = -4627

$ws.Range("H106").Value = 8706.333000000001
$ws.Range("I106").Value = 4434.8237
$ws.Range("J106").Value = 14292.154
$ws.Range("K106").Value = 4434.8237
$ws.Range("L106").Value = 14292.154
$ws.Range("M106").Value = -3803.8237
$ws.Range("N106").Value = -15554.154

$ws.Range("H116").Value = 13326.066
$ws.Range("J116").Value = 14241.714
$ws.Range("L116").Value = 14241.714
$ws.Range("N116").Value = -21125.714

$ws.Range("H132").Value = 1455.5641
$ws.Range("I132").Value = 1468.1892
$ws.Range("J132").Value = 1222
$ws.Range("K132").Value = 4404.5676
$ws.Range("L132").Value = 3666
$ws.Range("M132").Value = -1874.5676
$ws.Range("N132").Value = -8726

$ws.Range("H137").Value = 14927651
$ws.Range("I137").Value = 27779528
$ws.Range("J137").Value = 2891.2258
$ws.Range("K137").Value = 83338584
$ws.Range("L137").Value = 8673.6774
$ws.Range("M137").Value = -83336034
$ws.Range("N137").Value = -13773.6774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19234358
$ws.Range("I32").Value = 20835910
$ws.Range("K32").Value = 20835910
$ws.Range("M32").Value = -20835623

$ws.Range("H102").Value = 3999.5
$ws.Range("I102").Value = 4008.4546
$ws.Range("J102").Value = 3966.6667
$ws.Range("K102").Value = 4008.4546
$ws.Range("L102").Value = 3966.6667
$ws.Range("M102").Value = -2386.4546
$ws.Range("N102").Value = -7210.6667

$ws.Range("H132").Value = 2364.0232
$ws.Range("I132").Value = 1913.8684
$ws.Range("K132").Value = 5741.6052
$ws.Range("M132").Value = -3211.6052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2773.1304
$ws.Range("I134").Value = 1227.4
$ws.Range("K134").Value = 3682.2
$ws.Range("M134").Value = -1147.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41384.332
$ws.Range("I31").Value = 1046.4667
$ws.Range("J31").Value = 91806.664
$ws.Range("K31").Value = 1046.4667
$ws.Range("L31").Value = 91806.664
$ws.Range("M31").Value = -751.4666999999999
$ws.Range("N31").Value = -92396.664

$ws.Range("H34").Value = 41384.332
$ws.Range("I34").Value = 1046.4667
$ws.Range("J34").Value = 91806.664
$ws.Range("K34").Value = 1046.4667
$ws.Range("L34").Value = 91806.664
$ws.Range("M34").Value = -844.4666999999999
$ws.Range("N34").Value = -92210.664

$ws.Range("H39").Value = 7687.5
$ws.Range("I39").Value = 5375.5
$ws.Range("J39").Value = 9999.5
$ws.Range("K39").Value = 5375.5
$ws.Range("L39").Value = 9999.5
$ws.Range("M39").Value = -4984.5
$ws.Range("N39").Value = -10781.5

$ws.Range("H49").Value = 7687.5
$ws.Range("I49").Value = 5375.5
$ws.Range("J49").Value = 9999.5
$ws.Range("K49").Value = 5375.5
$ws.Range("L49").Value = 9999.5
$ws.Range("M49").Value = -5193.5
$ws.Range("N49").Value = -10363.5

$ws.Range("H105").Value = 6723.357
$ws.Range("I105").Value = 6168
$ws.Range("K105").Value = 6168
$ws.Range("M105").Value = -4421

$ws.Range("H107").Value = 1015.2
$ws.Range("I107").Value = 1020.8889
$ws.Range("K107").Value = 1020.8889
$ws.Range("M107").Value = 899.1111

$ws.Range("H131").Value = 43054.332
$ws.Range("I131").Value = 26000
$ws.Range("K131").Value = 26000
$ws.Range("M131").Value = -20960

$ws.Range("H132").Value = 3200.147
$ws.Range("I132").Value = 2135.8076
$ws.Range("K132").Value = 6407.4228
$ws.Range("M132").Value = -3877.4228

$ws.Range("H134").Value = 3797.4773
$ws.Range("I134").Value = 2458.3447
$ws.Range("J134").Value = 6386.467
$ws.Range("K134").Value = 7375.034100000001
$ws.Range("L134").Value = 19159.401
$ws.Range("M134").Value = -4840.034100000001
$ws.Range("N134").Value = -24229.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3601
$ws.Range("I70").Value = 3601
$ws.Range("K70").Value = 10803
$ws.Range("M70").Value = -10488

$ws.Range("H73").Value = 3601
$ws.Range("I73").Value = 3601
$ws.Range("K73").Value = 10803
$ws.Range("M73").Value = -9711

$ws.Range("H75").Value = 83339120
$ws.Range("I75").Value = 200002100
$ws.Range("J75").Value = 8421.143
$ws.Range("K75").Value = 600006300
$ws.Range("L75").Value = 25263.429
$ws.Range("M75").Value = -600005302
$ws.Range("N75").Value = -27259.429

$ws.Range("H78").Value = 83339120
$ws.Range("I78").Value = 200002100
$ws.Range("J78").Value = 8421.143
$ws.Range("K78").Value = 1800018900
$ws.Range("L78").Value = 75790.287
$ws.Range("M78").Value = -1800013908
$ws.Range("N78").Value = -85774.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888

$ws.Range("H113").Value = 1215
$ws.Range("I113").Value = 1217.5
$ws.Range("K113").Value = 1217.5
$ws.Range("M113").Value = 952.5

$ws.Range("H122").Value = 12544.546
$ws.Range("I122").Value = 16497.857
$ws.Range("K122").Value = 49493.571
$ws.Range("M122").Value = -47043.571

$ws.Range("H132").Value = 241655.3
$ws.Range("I132").Value = 287989.78
$ws.Range("K132").Value = 863969.3400000001
$ws.Range("M132").Value = -861439.3400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2055.3157
$ws.Range("I22").Value = 771.6667
$ws.Range("J22").Value = 2647.7693
$ws.Range("K22").Value = 771.6667
$ws.Range("L22").Value = 2647.7693
$ws.Range("M22").Value = -476.6667
$ws.Range("N22").Value = -3237.7693

$ws.Range("H27").Value = 2055.3157
$ws.Range("I27").Value = 771.6667
$ws.Range("J27").Value = 2647.7693
$ws.Range("K27").Value = 771.6667
$ws.Range("L27").Value = 2647.7693
$ws.Range("M27").Value = -664.6667
$ws.Range("N27").Value = -2861.7693

$ws.Range("H122").Value = 4818.5884
$ws.Range("I122").Value = 3400.1667
$ws.Range("J122").Value = 8222.799999999999
$ws.Range("K122").Value = 10200.5001
$ws.Range("L122").Value = 24668.4
$ws.Range("M122").Value = -7750.500100000001
$ws.Range("N122").Value = -29568.4

$ws.Range("H132").Value = 4394.0293
$ws.Range("I132").Value = 2400.3
$ws.Range("J132").Value = 7242.2144
$ws.Range("K132").Value = 7200.900000000001
$ws.Range("L132").Value = 21726.6432
$ws.Range("M132").Value = -4670.900000000001
$ws.Range("N132").Value = -26786.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3817.373
$ws.Range("I132").Value = 1488.3019
$ws.Range("J132").Value = 12634.571
$ws.Range("K132").Value = 4464.905699999999
$ws.Range("L132").Value = 37903.713
$ws.Range("M132").Value = -1934.905699999999
$ws.Range("N132").Value = -42963.713
